$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "LP1912" (index 1): date/time columns are B..G, new rows 1467-1482
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 01/01/2026 11:41:14"
$ws1.Range("A3").Value = "Total filas: 1481"

$sheet1Rows = @(
    @(1467, "11:41:03", "11:45", "215B_EL PATO", 4, "LP1912", "01/01/2026"),
    @(1468, "11:41:03", "11:52", "23_HERNANDEZ", 11, "LP1912", "01/01/2026"),
    @(1469, "11:41:03", "11:53", "225_GOMEZ", 12, "LP1912", "01/01/2026"),
    @(1470, "11:41:03", "12:01", "17_ROMERO", 20, "LP1912", "01/01/2026"),
    @(1471, "11:41:03", "12:06", "11_ETCHEVERRY", 25, "LP1912", "01/01/2026"),
    @(1472, "11:41:03", "12:16", "215C_EL PATO", 35, "LP1912", "01/01/2026"),
    @(1473, "11:41:03", "12:31", "15_ABASTO", 50, "LP1912", "01/01/2026"),
    @(1474, "11:41:03", "12:31", "23_HERNANDEZ", 50, "LP1912", "01/01/2026"),
    @(1475, "11:41:03", "12:33", "14_ABASTO", 52, "LP1912", "01/01/2026"),
    @(1476, "11:41:03", "12:37", "27_EL RETIRO", 56, "LP1912", "01/01/2026"),
    @(1477, "11:41:03", "12:46", "15X38_ABASTO", 65, "LP1912", "01/01/2026"),
    @(1478, "11:41:03", "12:48", "16_SANTA ANA", 67, "LP1912", "01/01/2026"),
    @(1479, "11:41:03", "13:00", "16_SANTA ANA", 79, "LP1912", "01/01/2026"),
    @(1480, "11:41:03", "13:04", "215C_EL PATO", 83, "LP1912", "01/01/2026"),
    @(1481, "11:41:03", "13:13", "16_SANTA ANA", 92, "LP1912", "01/01/2026"),
    @(1482, "11:41:03", "13:17", "10_OLMOS", 96, "LP1912", "01/01/2026")
)

# Pre-format the new "Fecha" cells as text so the "dd/mm/yyyy"-looking
# string is stored verbatim (as every other row in this column already is)
# instead of being auto-converted into a date serial number.
$ws1.Range("G1467:G1482").NumberFormat = "@"

foreach ($row in $sheet1Rows) {
    $r = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $ws1.Cells.Item($r, 4).Value = $row[3]
    $ws1.Cells.Item($r, 5).Value = $row[4]
    $ws1.Cells.Item($r, 6).Value = $row[5]
    $ws1.Cells.Item($r, 7).Value = $row[6]
}

# ---------------------------------------------------------------------------
# Sheet "LP1912-215" (index 2): date/time columns are B..G, new rows 105-107
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 01/01/2026 11:41:14"
$ws2.Range("A3").Value = "Total filas: 106"

$sheet2Rows = @(
    @(105, "01/01/2026", "11:41:03", "11:45", "215B_EL PATO", 4, "LP1912"),
    @(106, "01/01/2026", "11:41:03", "12:16", "215C_EL PATO", 35, "LP1912"),
    @(107, "01/01/2026", "11:41:03", "13:04", "215C_EL PATO", 83, "LP1912")
)

$ws2.Range("B105:B107").NumberFormat = "@"

foreach ($row in $sheet2Rows) {
    $r = $row[0]
    $ws2.Cells.Item($r, 2).Value = $row[1]
    $ws2.Cells.Item($r, 3).Value = $row[2]
    $ws2.Cells.Item($r, 4).Value = $row[3]
    $ws2.Cells.Item($r, 5).Value = $row[4]
    $ws2.Cells.Item($r, 6).Value = $row[5]
    $ws2.Cells.Item($r, 7).Value = $row[6]
}

# ---------------------------------------------------------------------------
# Sheet "6203-6173" (index 3): date/time columns are B..G, new row 184
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 01/01/2026 11:41:14"
$ws3.Range("A3").Value = "Total filas: 183"

$sheet3Rows = @(
    ,@(184, "01/01/2026", "11:41:09", "13:06", "215C_LA PLATA", 85, "L6203")
)

$ws3.Range("B184:B184").NumberFormat = "@"

foreach ($row in $sheet3Rows) {
    $r = $row[0]
    $ws3.Cells.Item($r, 2).Value = $row[1]
    $ws3.Cells.Item($r, 3).Value = $row[2]
    $ws3.Cells.Item($r, 4).Value = $row[3]
    $ws3.Cells.Item($r, 5).Value = $row[4]
    $ws3.Cells.Item($r, 6).Value = $row[5]
    $ws3.Cells.Item($r, 7).Value = $row[6]
}

Write-Output "edit complete"
